$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Abhishek Pathak, EC3219): Present=3, Total Classes=3, Attendance %=100, Attendance Range="91-100%"
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = "91-100%"

# Row 3 (Shubham Pitekar, EC3230): Present=2, Total Classes=2 (Attendance % and Range unchanged)
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 2
